$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2018-10", 101, 97.8, 99.2),
    @("2018-11", 101.2, 97.8, 99.59999999999999),
    @("2018-12", 100.8, 97.7, 99.8),
    @("2018-01", 102.6615, 98.8289, 101.1681),
    @("2018-02", 102.1, 98.7, 101.4),
    @("2018-03", 102.1, 98.40000000000001, 101.6),
    @("2018-04", 101.3, 98.7, 101),
    @("2018-05", 101.4, 98.2, 101.2),
    @("2018-06", 101.5, 98, 101),
    @("2018-07", 101, 97.90000000000001, 100.6),
    @("2018-08", 101, 98.09999999999999, 99.90000000000001),
    @("2018-09", 101.1, 98.5, 99.3),
    @("2019-10", 101.1, 98.59999999999999, 100.8),
    @("2019-11", 100.9, 98.40000000000001, 100.6),
    @("2019-12", 101.1, 98.59999999999999, 100.5),
    @("2019-01", 100.9, 97.90000000000001, 99.7),
    @("2019-02", 100.9, 97.90000000000001, 99.90000000000001),
    @("2019-03", 100.7, 98, 100),
    @("2019-04", 101.1, 98.2, 101),
    @("2019-05", 101.2, 98.7, 101.2),
    @("2019-06", 101, 98.8, 100.9),
    @("2019-07", 101, 98.2, 100.7),
    @("2019-08", 101.3, 98, 100.5),
    @("2019-09", 101, 97.90000000000001, 100.8),
    @("2020-10", 100.8, 97.8, 98.5),
    @("2020-11", 100.6, 97.90000000000001, 98.5),
    @("2020-12", 100.7, 98.09999999999999, 99.3),
    @("2020-01", 100.6, 98.40000000000001, 100.8),
    @("2020-02", 100.6, 98, 100.6),
    @("2020-03", 100.5, 97.59999999999999, 99.90000000000001),
    @("2020-04", 100.3, 96.5, 99.59999999999999),
    @("2020-05", 100.3, 96.3, 99),
    @("2020-06", 100.4, 96.40000000000001, 99.09999999999999),
    @("2020-07", 100.6, 97, 98.8),
    @("2020-08", 100.2, 97.3, 98.7),
    @("2020-09", 100.5, 97.7, 98.2),
    @("2021-10", 104.3, 100.6, 100.6),
    @("2021-11", 107.3, 101.5, 103.4),
    @("2021-12", 109.1, 103.6, 106.5),
    @("2021-01", 101.7, 97.8, 98.59999999999999),
    @("2021-02", 101.8, 98.5, 98.59999999999999),
    @("2021-03", 101.9, 99.09999999999999, 99),
    @("2021-04", 101.4, 99.8, 98.8),
    @("2021-05", 101.2, 99.90000000000001, 98.7),
    @("2021-06", 101.5, 99.7, 98.59999999999999),
    @("2021-07", 102.1, 100.1, 99),
    @("2021-08", 102.4, 100, 98.90000000000001),
    @("2021-09", 103.1, 100.4, 99.7),
    @("2022-10", 104.2, 108.7, 110.4),
    @("2022-11", 101, 107.7, 108.2),
    @("2022-12", 100.1, 105.5, 105),
    @("2022-01", 108.6, 107.3, 108.2),
    @("2022-02", 107.5, 107.8, 109.7),
    @("2022-03", 106.3, 108, 110),
    @("2022-04", 107, 108.2, 111.3),
    @("2022-05", 107.7, 108.2, 112.3),
    @("2022-06", 106.8, 107.8, 112.8),
    @("2022-07", 106.6, 106.7, 111.4),
    @("2022-08", 105.9, 107.6, 112),
    @("2022-09", 105.6, 108.2, 113.6),
    @("2023-01", 101.3, 102.7, 103.4),
    @("2023-02", 101.5, 102.9, 103),
    @("2023-03", 101.8, 102, 102.5),
    @("2023-04", 100.5, 101.4, 102),
    @("2023-05", 99.7, 101.5, 102),
    @("2023-06", 100.3, 101.3, 101.3),
    @("2023-07", 99.5, 101.8, 102.1)
)

$rowCount = $data.Count

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = [double]$row[1]
    $ws.Cells.Item($r, 3).Value = [double]$row[2]
    $ws.Cells.Item($r, 4).Value = [double]$row[3]
}

# The original sheet only had formatting defined through row 49 (A1:D49).
# Newly-added rows beyond that need the same per-cell formatting (bold/
# bordered/centered date cells in column A) as the pre-existing rows, so
# copy the formatting from an existing fully-formatted row down across the
# freshly added rows.
$lastOriginalRow = 49
$lastNewRow = $rowCount + 1
if ($lastNewRow -gt $lastOriginalRow) {
    $fmtSrc = $ws.Range("A2:D2")
    $fmtSrc.Copy()
    $fmtDst = $ws.Range("A" + ($lastOriginalRow + 1) + ":D" + $lastNewRow)
    $fmtDst.PasteSpecial(-4122)
}
